# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worksheet lists overdue-payment periods per worker for NIT 9002496029.
# This update replaces the worker/period data block with a new dataset
# (3 workers instead of 6) and refreshes the summary totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Remove 3 rows from the existing data block (rows 16-29 all share the
#    same "interior row" style; row 30 carries the unique "closing row"
#    bottom-border style). Deleting 3 interior rows shrinks the 15-row
#    block (16-30, 6 workers) down to 12 rows (16-27, 3 workers) while
#    automatically sliding row 30's distinct style up to become the new
#    row 27, and shifting every row below (footer, merges) up by 3.
# ---------------------------------------------------------------------
$ws.Rows("16:18").Delete()

# ---------------------------------------------------------------------
# 2. Write the new worker / period rows into B16:G27.
# ---------------------------------------------------------------------
$data = @(
    @("CC", "1143378856", "CRISTIAN PITALUA BATISTA", "2012", 7022, 877803),
    @("CC", "1143345476", "JUAN CARLOS CALVO CASTILLO", "2107", 33918, 908526),
    @("CC", "1143345476", "JUAN CARLOS CALVO CASTILLO", "2108", 36341, 908526),
    @("CC", "1143345476", "JUAN CARLOS CALVO CASTILLO", "2109", 36341, 908526),
    @("CC", "1143345476", "JUAN CARLOS CALVO CASTILLO", "2110", 36341, 908526),
    @("CC", "1143345476", "JUAN CARLOS CALVO CASTILLO", "2111", 30284, 908526),
    @("CC", "1048442306", "JAILER ENRIQUE BOHORQUEZ GOMEZ", "2501", 39858, 1423500),
    @("CC", "1048442306", "JAILER ENRIQUE BOHORQUEZ GOMEZ", "2502", 56940, 1423500),
    @("CC", "1048442306", "JAILER ENRIQUE BOHORQUEZ GOMEZ", "2503", 56940, 1423500),
    @("CC", "1048442306", "JAILER ENRIQUE BOHORQUEZ GOMEZ", "2504", 56940, 1423500),
    @("CC", "1048442306", "JAILER ENRIQUE BOHORQUEZ GOMEZ", "2505", 56940, 1423500),
    @("CC", "1048442306", "JAILER ENRIQUE BOHORQUEZ GOMEZ", "2506", 56940, 1423500)
)

$row = 16
foreach ($rec in $data) {
    $ws.Range("B$row").Value = $rec[0]
    $ws.Range("C$row").Value = $rec[1]
    $ws.Range("D$row").Value = $rec[2]
    $ws.Range("E$row").Value = $rec[3]
    $ws.Range("F$row").Value = $rec[4]
    $ws.Range("G$row").Value = $rec[5]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 3. Refresh the summary block above the table.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 504805   # VALOR MORA total
$ws.Range("C13").Value = 3        # Cant. Trabajadores
$ws.Range("F13").Value = 12       # Cant. Periodos
